$d = $word.ActiveDocument

# --- 1. Merge the "Ejemplos" + ":" runs into a single run ("Ejemplos:") ---
$d.Content.Find.Execute("Ejemplos:", $true, $false, $false, $false, $false, $true, 1, $false, "Ejemplos:", 2) | Out-Null

# --- 2. Remove the _GoBack bookmark from its old location ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 3. Add a new paragraph at the end of the document with "Prueba commit" ---
$sel = $word.Selection
$sel.EndKey(6, 0)
$sel.TypeParagraph()

$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.MoveEnd(1, -1)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>Prueba commit</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

# InsertXML leaves a vestigial empty trailing paragraph behind; remove it along with
# the paragraph mark that separates it from our new "Prueba commit" paragraph.
$endPos = $d.Content.End
$delRange = $d.Range($endPos - 2, $endPos - 1)
$delRange.Delete()

# --- 4. Re-create the _GoBack bookmark right after "Prueba commit" ---
# Type a disambiguating trailing character so the collapsed bookmark range is placed
# *after* the run rather than wrapping it, then remove that character again.
$sel2 = $word.Selection
$sel2.EndKey(6, 0)
$sel2.TypeText("X")

$markPos = $d.Content.End - 1
$bmRange = $d.Range($markPos, $markPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$cleanupPos = $d.Content.End - 1
$cleanupRange = $d.Range($cleanupPos, $cleanupPos + 1)
$cleanupRange.Delete()
